# Update the EC (Estado de Cuenta) worksheet:
#  - Replace the worker dataset (was 4 workers across several periods,
#    now a single worker across periods 2505-2508)
#  - Update the summary figures (VALOR MORA, Cant. Trabajadores)
#  - Update the Salario Basico column for the first two data rows
#  - Remove now-unused data rows (the row count shrank from 11 to 4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Delete the data rows that are no longer needed -----------------------
# Rows 16-18 stay (their position/style is reused for the remaining data),
# but the old rows 19-25 (7 rows) must go so that the table's last styled
# row (formerly row 26, with the heavier bottom border) slides up into the
# row-19 slot, and the signature rows (formerly 31/32) slide up into 24/25.
$ws.Rows("19:25").Delete()

# --- Summary block ----------------------------------------------------------
$ws.Range("E11").Value = 227760      # VALOR MORA
$ws.Range("C13").Value = 1           # Cant. Trabajadores

# --- Data rows (single worker YESSIKA YIBETH PEREZ ALVAREZ, doc 1047502073) -
$ws.Range("C16").Value = "1047502073"
$ws.Range("D16").Value = "YESSIKA YIBETH PEREZ ALVAREZ"
$ws.Range("E16").Value = "2505"
$ws.Range("G16").Value = 1423500

$ws.Range("C17").Value = "1047502073"
$ws.Range("D17").Value = "YESSIKA YIBETH PEREZ ALVAREZ"
$ws.Range("E17").Value = "2506"
$ws.Range("G17").Value = 1423500

$ws.Range("C18").Value = "1047502073"
$ws.Range("D18").Value = "YESSIKA YIBETH PEREZ ALVAREZ"
$ws.Range("E18").Value = "2507"
$ws.Range("G18").Value = 1423500

$ws.Range("C19").Value = "1047502073"
$ws.Range("D19").Value = "YESSIKA YIBETH PEREZ ALVAREZ"
$ws.Range("E19").Value = "2508"

# --- Column D was sized for a shorter name now; narrow it a touch ----------
$ws.Columns("D").ColumnWidth = 30
